$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (item 1): BAV21A0F1C (ball valve) -> CKV413C3J2G (swing check valve)
$ws.Cells.Item(2, 2).Value = "CKV413C3J2G"
$ws.Cells.Item(2, 3).Value = "SWING CHECK VALVE FL, MSS SP-136, A536 Gr 65-45-12, CL 125, INST HORIZ/VERT, FF, B16.1, BOLTED COVER, NON METALLIC FLAT GASKET EPDM, SEATS BRONZE, DISC DUCTILE IRON EPDM ENCAPSULATED, SS304 PIN,UL LISTED/FM APPROVED"
$ws.Cells.Item(2, 4).Value = "6,00"
$ws.Cells.Item(2, 5).Value = "1,00"
$ws.Cells.Item(2, 6).Value = "e.a"
$ws.Cells.Item(2, 7).Value = "-"

# Row 3 (item 2): BAV21A0F1C (ball valve) -> GAV413C3J2G (gate valve)
$ws.Cells.Item(3, 2).Value = "GAV413C3J2G"
$ws.Cells.Item(3, 3).Value = "GATE VALVE FL, MSS SP-128 TYPE II, A536 Gr 65-45-12, CL 125, FF, B16.1, BB, NON METALLIC FLAT GASKET EPDM, PKG EPDM, SOLID WEDGE, WEDGE DUCTILE IRON EPDM ENCAPSULATED, STEM BRONZE, OS&Y/RSNRO, HO, UL LISTED/FM APPROVED"
$ws.Cells.Item(3, 4).Value = "6,00"
$ws.Cells.Item(3, 5).Value = "9,00"
$ws.Cells.Item(3, 6).Value = "e.a"
$ws.Cells.Item(3, 7).Value = "-"

# Row 4 (item 3): BAV21A0F1C (ball valve) -> GAV413C3J2G (gate valve)
$ws.Cells.Item(4, 2).Value = "GAV413C3J2G"
$ws.Cells.Item(4, 3).Value = "GATE VALVE FL, MSS SP-128 TYPE II, A536 Gr 65-45-12, CL 125, FF, B16.1, BB, NON METALLIC FLAT GASKET EPDM, PKG EPDM, SOLID WEDGE, WEDGE DUCTILE IRON EPDM ENCAPSULATED, STEM BRONZE, OS&Y/RSNRO, HO, UL LISTED/FM APPROVED"
$ws.Cells.Item(4, 4).Value = "8,00"
$ws.Cells.Item(4, 5).Value = "1,00"
$ws.Cells.Item(4, 6).Value = "e.a"
$ws.Cells.Item(4, 7).Value = "-"

# Row 5 (item 4): BAV21A0F1C (ball valve) -> GAV41H4J2G (gate valve, threaded)
$ws.Cells.Item(5, 2).Value = "GAV41H4J2G"
$ws.Cells.Item(5, 3).Value = "GATE VALVE THD, MSS SP-80, B62 UNS C83600, CL 150, NPTF, B1.20.1, SCREWED BONNET, PKG LUBRICATED FIBER/GRAPH, BRONZE STEM, SEAT RINGS&DISC, S, SOLID WEDGE, STEM OS&Y/RSNRO, HO"
$ws.Cells.Item(5, 4).Value = "0,5"
$ws.Cells.Item(5, 5).Value = "1,00"
$ws.Cells.Item(5, 6).Value = "e.a"
$ws.Cells.Item(5, 7).Value = "-"

# Row 6 (item 5): BAV24G0I1C (ball valve) -> GLV415J4J2G (angle hose valve)
$ws.Cells.Item(6, 2).Value = "GLV415J4J2G"
$ws.Cells.Item(6, 3).Value = "ANGLE HOSE VALVE THD, MSS SP-80, B62 UNS C83600, 300 PSI CWP, FNPT X NH W/CAP AND CHAIN, FNPT B1.20.1 AND NH NFPA 1963, SCREW-IN BONNET, RENEWABLE DISC, PKG NON ASBESTOS, DISC&STEM BRONZE, RISING STEM, HO, UL LISTED/FM APPROVED"
$ws.Cells.Item(6, 4).Value = "2,5"
$ws.Cells.Item(6, 5).Value = "4,00"
$ws.Cells.Item(6, 6).Value = "e.a"
$ws.Cells.Item(6, 7).Value = "-"

# Remove rows 7-11 (items 6-10: CKV21A0B2B, CLV24F0B2B, MFV21A0I2I x3)
$ws.Range("A7:G11").EntireRow.Delete()
